$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = 3
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0.003
$ws.Cells.Item(28, 5).Value = "Regular"
$ws.Cells.Item(28, 9).Value = "<function relu at 0x10f4d69d8>"
$ws.Cells.Item(28, 10).Value = 0.9516000151634216
$ws.Cells.Item(28, 11).Value = 0.04560000076889992
$ws.Cells.Item(28, 12).Value = 0.003700000001117587
$ws.Cells.Item(28, 13).Value = 0.1825118958950043
$ws.Cells.Item(28, 14).Value = 6.906796932220459
$ws.Cells.Item(28, 15).Value = 0.04560000076889992
$ws.Cells.Item(28, 16).Value = "logs/results_282.log"
$ws.Cells.Item(28, 17).Value = "weights/model_282.ckpt"
$ws.Cells.Item(28, 18).Value = "tb/282"
$ws.Cells.Item(28, 19).Value = "(6.9546156, 7.3834124, 9.092276, 9.433221, 9.5498905, 11.14911, 9.904368)"
$ws.Cells.Item(28, 20).Value = "(139.42224, 9.018682, 9.271418, 8.975029, 7.9992733, 7.393931, 7.053868, 10.015819)"

# Row 29
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 3
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 0.003
$ws.Cells.Item(29, 5).Value = "Regular"
$ws.Cells.Item(29, 9).Value = "<function relu at 0x1100289d8>"
$ws.Cells.Item(29, 10).Value = 0.9430999755859375
$ws.Cells.Item(29, 11).Value = 0.06279999762773514
$ws.Cells.Item(29, 12).Value = 0.006200000178068876
$ws.Cells.Item(29, 13).Value = 0.2192680686712265
$ws.Cells.Item(29, 14).Value = 7.151318073272705
$ws.Cells.Item(29, 15).Value = 0.06279999762773514
$ws.Cells.Item(29, 16).Value = "logs/results_285.log"
$ws.Cells.Item(29, 17).Value = "weights/model_285.ckpt"
$ws.Cells.Item(29, 18).Value = "tb/285"
$ws.Cells.Item(29, 19).Value = "(7.0936155, 7.6124697, 8.267413, 8.219525, 11.057663, 10.19839, 9.6345)"
$ws.Cells.Item(29, 20).Value = "(141.34113, 8.731318, 9.105043, 8.344593, 9.322138, 7.830576, 7.469233, 9.507704)"

# Row 30
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 3
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0.003
$ws.Cells.Item(30, 5).Value = "Regular"
$ws.Cells.Item(30, 9).Value = "<function relu at 0x116aa89d8>"
$ws.Cells.Item(30, 10).Value = 0.9474999904632568
$ws.Cells.Item(30, 11).Value = 0.07349999994039536
$ws.Cells.Item(30, 12).Value = 0.03050000034272671
$ws.Cells.Item(30, 13).Value = 0.2007102072238922
$ws.Cells.Item(30, 14).Value = 6.25114631652832
$ws.Cells.Item(30, 15).Value = 0.07349999994039536
$ws.Cells.Item(30, 16).Value = "logs/results_305.log"
$ws.Cells.Item(30, 17).Value = "weights/model_305.ckpt"
$ws.Cells.Item(30, 18).Value = "tb/305"
$ws.Cells.Item(30, 19).Value = "(6.9461164, 7.606389, 7.896417, 8.789286, 9.165759, 8.689637, 8.321884)"
$ws.Cells.Item(30, 20).Value = "(138.02159, 8.655811, 9.216804, 9.208384, 7.499287, 7.888038, 7.191258, 13.024145)"
